# Animals.docx word-list edit
# Rebuild the word list paragraphs from scratch (this naturally drops the
# legacy <w:proofErr gramStart/gramEnd> markup that Word had stamped on
# every paragraph - fresh content never carries that), expand the list to
# 40 animal words, re-anchor the _GoBack bookmark so it spans the whole
# list (start right before the first word, end right after the last),
# and make sure the document ends with two empty paragraphs.

$d = $word.ActiveDocument

# ---- 1. wipe all existing body content down to a single empty paragraph
$guard = 0
while ($d.Content.End -gt 1 -and $guard -lt 500) {
    $d.Content.Delete()
    $guard = $guard + 1
}

# ---- 2. the 40-word list (order matters - matches the target word count
#         called out in the commit message: "40 words in each word list")
$words = @(
    "monkey","reptile","pterodactyl ","gerbil ","parakeet ","dolphin",
    "reindeer","cockatoo","geese","sparrow","clownfish","lemur",
    "mockingbird","amphibian","mosquito","sloth","spider","lizard",
    "woodpecker","feline","dragon","unicorn","pegasus","dinosaur",
    "triceratops","gorilla","fly","goldfish","hamster","sheep","beetle",
    "buffalo","mallard","grasshopper","armadillo","antelope","alligator",
    "alpaca","llama","elephant"
)

# ---- 3. write the words out, one per paragraph
for ($i = 0; $i -lt $words.Count; $i++) {
    if ($i -gt 0) {
        $p = $d.Content
        $p.Collapse(0)
        $p.InsertParagraphAfter()
    }
    $r = $d.Content
    $r.Collapse(0)
    $r.InsertAfter($words[$i])
}

# ---- 4. one more trailing empty paragraph (target ends with two <w:p/>)
$p = $d.Content
$p.Collapse(0)
$p.InsertParagraphAfter()

# ---- 5. re-anchor the _GoBack bookmark to span the whole word list:
#         from right before "monkey" through right after "elephant".
$first = $d.Paragraphs.Item(1).Range
$wordCount = $words.Count
$lastWordParagraph = $d.Paragraphs.Item($wordCount)
$span = $d.Range($first.Start, $lastWordParagraph.Range.End)
$d.Bookmarks.Add("_GoBack", $span)
